# Append 5 new rows (7-11) to the worksheet. These duplicate earlier
# match rows (re-scraped data appended in a different order), matching
# the json-file-driven re-scrape described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(" Abu Dhabi", " October 25 2020", "Royals won by 8 wickets (with 10 balls remaining)", "Mumbai Indians", "Rajasthan Royals", "Saurabh Tiwary ", "34", "25", "4", "1", "136.00"),
    @(" Abu Dhabi", " September 19 2020", "Super Kings won by 5 wickets (with 4 balls remaining)", "Mumbai Indians", "Chennai Super Kings", "Saurabh Tiwary ", "42", "31", "3", "1", "135.48"),
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Mumbai Indians", "Royal Challengers Bangalore", "Saurabh Tiwary ", "5", "8", "0", "0", "62.50"),
    @(" Abu Dhabi", " September 23 2020", "Mumbai won by 49 runs", "Mumbai Indians", "Kolkata Knight Riders", "Saurabh Tiwary ", "21", "13", "1", "1", "161.53"),
    @(" Sharjah", " November 03 2020", "Sunrisers won by 10 wickets (with 17 balls remaining)", "Mumbai Indians", "Sunrisers Hyderabad", "Saurabh Tiwary ", "1", "3", "0", "0", "33.33")
)

$startRow = 7
$endRow = $startRow + $newRows.Count - 1

# Keep all the new values as text (the sheet stores every column,
# including the numeric-looking ones, as text -- see the existing
# ignoredErrors/numberStoredAsText on the sheet).
$ws.Range("A$startRow`:K$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
